$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "להלה"
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "Love Gift Card"
$ws.Range("E2").Value = "שובר בסך 100 ש""ח"
$ws.Range("F2").Value = "'2023-10-31"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "מבצע"
$ws.Range("H2").Value = "פעיל"
